# NIT-9005910705.xlsx - "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# Rebuilds the worker debt table with the new dataset (6 workers / 15 period rows)
# and updates the summary totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room: the table grows from 10 data rows (16-25) to 15 data rows
#    (16-30). Insert 5 rows just above the old "closing" row (25) so that row
#    25 (with its special bottom-border style) ends up at row 30, and the
#    footer block (old rows 30-31) ends up at rows 35-36.
# ---------------------------------------------------------------------------
$ws.Rows("25:29").Insert(-4121)   # xlShiftDown

# Copy the formatting (borders/number formats) of a normal interior data row
# (row 24) down into the 5 newly inserted blank rows.
$ws.Range("B24:J24").Copy()
$ws.Range("B25:J29").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Fill the table with the new dataset.
#    Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#             E=Periodo Mora, F=Valor Mora, G=Salario Basico
# ---------------------------------------------------------------------------
$data = @(
    @{Row=16; C="73196146";    D="ANDRES MIGUEL VILLADIEGO OTERO";    E="2307"; F=10667; G=1000000},
    @{Row=17; C="73196146";    D="ANDRES MIGUEL VILLADIEGO OTERO";    E="2306"; F=40000; G=1000000},
    @{Row=18; C="73196146";    D="ANDRES MIGUEL VILLADIEGO OTERO";    E="2305"; F=40000; G=1000000},
    @{Row=19; C="73196146";    D="ANDRES MIGUEL VILLADIEGO OTERO";    E="2304"; F=40000; G=1000000},
    @{Row=20; C="73196146";    D="ANDRES MIGUEL VILLADIEGO OTERO";    E="2303"; F=40000; G=1000000},
    @{Row=21; C="73196146";    D="ANDRES MIGUEL VILLADIEGO OTERO";    E="2302"; F=40000; G=1000000},
    @{Row=22; C="1098750705";  D="EDGAR MAURICIO VALENCIA CARRILLO";  E="2507"; F=28470; G=1423500},
    @{Row=23; C="1047367540";  D="WILLIAN ENRIQUE ARROYO OROZCO";     E="2302"; F=40000; G=1000000},
    @{Row=24; C="1047367540";  D="WILLIAN ENRIQUE ARROYO OROZCO";     E="2301"; F=40000; G=1000000},
    @{Row=25; C="1047367540";  D="WILLIAN ENRIQUE ARROYO OROZCO";     E="2211"; F=40000; G=1000000},
    @{Row=26; C="1047367540";  D="WILLIAN ENRIQUE ARROYO OROZCO";     E="2209"; F=40000; G=1000000},
    @{Row=27; C="20310060";    D="LESGUI CUADRADO BANQUEZ";           E="2106"; F=18170; G=908526},
    @{Row=28; C="20337315";    D="RONELIS BERRIO ARIAS";              E="2112"; F=29073; G=908526},
    @{Row=29; C="91146070";    D="CHRISTIAN JOSE RAMIREZ ALARCON";    E="2207"; F=4000;  G=1000000},
    @{Row=30; C="91146070";    D="CHRISTIAN JOSE RAMIREZ ALARCON";    E="2206"; F=32000; G=1000000}
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = $rec.C
    $ws.Range("D$r").Value = $rec.D
    $ws.Range("E$r").Value = $rec.E
    $ws.Range("F$r").Value = $rec.F
    $ws.Range("G$r").Value = $rec.G
}

# ---------------------------------------------------------------------------
# 3. Update the summary block above the table.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 482380    # VALOR MORA (matches sum of column F)
$ws.Range("C13").Value = 6         # Cant. Trabajadores
$ws.Range("F13").Value = 14        # Cant. Periodos

# ---------------------------------------------------------------------------
# 4. Column D widened slightly to fit the longer new name
#    ("EDGAR MAURICIO VALENCIA CARRILLO").
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 35
